$d = $word.ActiveDocument

function Remove-TrailingUnderlineRuns($paraIndex, $anchorText) {
    # Finds $anchorText (the end of the "real" paragraph content) and deletes
    # everything between the end of that match and the end of the paragraph
    # (excluding the paragraph mark) -- this strips the extra underlined
    # "quiz word" runs that trail the real text in a paragraph.
    $para = $d.Paragraphs($paraIndex)
    $searchRange = $para.Range.Duplicate
    $found = $searchRange.Find.Execute($anchorText)
    if ($found) {
        $delStart = $searchRange.End
        $delEnd = $para.Range.End - 1
        if ($delEnd -gt $delStart) {
            $del = $d.Range($delStart, $delEnd)
            $del.Delete()
        }
    }
}

# --- Paragraph 2 (Fish body): merge "H. Rackham" break + leading quote into ":" ---
$d.Content.Find.Execute("H. Rackham" + [char]11 + [char]34, $true, $false, $false, $false, $false, $true, 1, $false, "H. Rackham:", 2) | Out-Null

# Remove the trailing quote mark right before the first break pair (end of Rackham quote)
$d.Content.Find.Execute("resultant pleasure?" + [char]34, $true, $false, $false, $false, $false, $true, 1, $false, "resultant pleasure?", 2) | Out-Null

# Remove the leading quote mark before the Cicero quote
$d.Content.Find.Execute([char]34 + "At vero eos", $true, $false, $false, $false, $false, $true, 1, $false, "At vero eos", 2) | Out-Null

# Remove the trailing quote mark after the Cicero quote
$d.Content.Find.Execute("asperiores repellat." + [char]34, $true, $false, $false, $false, $false, $true, 1, $false, "asperiores repellat.", 2) | Out-Null

# Strip the trailing underlined "quiz word" runs from paragraph 2
Remove-TrailingUnderlineRuns 2 "asperiores repellat."

# --- Paragraph 4 (Cheese body): strip wrapping quotes ---
$d.Content.Find.Execute([char]34 + "Sed ut perspiciatis", $true, $false, $false, $false, $false, $true, 1, $false, "Sed ut perspiciatis", 2) | Out-Null
$d.Content.Find.Execute("nulla pariatur?" + [char]34, $true, $false, $false, $false, $false, $true, 1, $false, "nulla pariatur.", 2) | Out-Null

# Strip the trailing underlined "quiz word" runs from paragraph 4
Remove-TrailingUnderlineRuns 4 "nulla pariatur."

# --- Paragraph 5 (Car heading): remove duplicated underlined "Car" run ---
Remove-TrailingUnderlineRuns 5 "Car"

# --- Paragraph 6 (Lorem ipsum body): strip wrapping quotes ---
$d.Content.Find.Execute([char]34 + "Lorem ipsum dolor", $true, $false, $false, $false, $false, $true, 1, $false, "Lorem ipsum dolor", 2) | Out-Null
$d.Content.Find.Execute("id est laborum." + [char]34, $true, $false, $false, $false, $false, $true, 1, $false, "id est laborum.", 2) | Out-Null

# Strip the trailing underlined "quiz word" runs from paragraph 6
Remove-TrailingUnderlineRuns 6 "id est laborum."
